$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 908.5833
$ws.Range("I2").Value = 963.7273
$ws.Range("J2").Value = 302
$ws.Range("K2").Value = 963.7273
$ws.Range("L2").Value = 302
$ws.Range("M2").Value = -850.7273
$ws.Range("N2").Value = -528
$ws.Range("H5").Value = 172.92308
$ws.Range("I5").Value = 80.75
$ws.Range("K5").Value = 80.75
$ws.Range("M5").Value = 34.25
$ws.Range("H9").Value = 6536282.5
$ws.Range("I9").Value = 462.5
$ws.Range("J9").Value = 22222250
$ws.Range("K9").Value = 462.5
$ws.Range("L9").Value = 22222250
$ws.Range("M9").Value = -293.5
$ws.Range("N9").Value = -22222588
$ws.Range("H19").Value = 1048.8889
$ws.Range("I19").Value = 500.33334
$ws.Range("J19").Value = 1323.1666
$ws.Range("K19").Value = 500.33334
$ws.Range("L19").Value = 1323.1666
$ws.Range("M19").Value = -325.33334
$ws.Range("N19").Value = -1673.1666
$ws.Range("H76").Value = 6469.4287
$ws.Range("I76").Value = 6296.3335
$ws.Range("K76").Value = 6296.3335
$ws.Range("M76").Value = -5981.3335
$ws.Range("H79").Value = 6469.4287
$ws.Range("I79").Value = 6296.3335
$ws.Range("K79").Value = 6296.3335
$ws.Range("M79").Value = -5204.3335
$ws.Range("H103").Value = 29673.438
$ws.Range("I103").Value = 46224.445
$ws.Range("J103").Value = 8393.571
$ws.Range("K103").Value = 138673.335
$ws.Range("L103").Value = 25180.713
$ws.Range("M103").Value = -138087.335
$ws.Range("N103").Value = -26352.713
$ws.Range("H129").Value = 1694.3684
$ws.Range("I129").Value = 1422.8334
$ws.Range("K129").Value = 4268.5002
$ws.Range("M129").Value = 731.4997999999996
$ws.Range("H132").Value = 4500.0835
$ws.Range("J132").Value = 7072.143
$ws.Range("L132").Value = 21216.429
$ws.Range("N132").Value = -26276.429
$ws.Range("H133").Value = 75551.664
$ws.Range("J133").Value = 75551.664
$ws.Range("L133").Value = 75551.664
$ws.Range("N133").Value = -85671.664
$ws.Range("H137").Value = 402145.53
$ws.Range("I137").Value = 688186.3
$ws.Range("K137").Value = 2064558.9
$ws.Range("M137").Value = -2062008.9
$ws.Range("H138").Value = 6470.4707
$ws.Range("J138").Value = 6852
$ws.Range("L138").Value = 20556
$ws.Range("N138").Value = -30836

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H13").Value = 4999.3335
$ws.Range("I13").Value = 4999
$ws.Range("K13").Value = 4999
$ws.Range("M13").Value = -4855
$ws.Range("H32").Value = 3462.35
$ws.Range("I32").Value = 3427.6843
$ws.Range("K32").Value = 3427.6843
$ws.Range("M32").Value = -3140.6843
$ws.Range("H34").Value = 211605
$ws.Range("J34").Value = 400000
$ws.Range("L34").Value = 400000
$ws.Range("N34").Value = -400542
$ws.Range("H43").Value = 19691.75
$ws.Range("I43").Value = 18867.334
$ws.Range("J43").Value = 20186.4
$ws.Range("K43").Value = 18867.334
$ws.Range("L43").Value = 20186.4
$ws.Range("M43").Value = -18554.334
$ws.Range("N43").Value = -20812.4
$ws.Range("H61").Value = 11353.094
$ws.Range("I61").Value = 13332.167
$ws.Range("K61").Value = 13332.167
$ws.Range("M61").Value = -13120.167
$ws.Range("H63").Value = 1073.5714
$ws.Range("I63").Value = 1103.4
$ws.Range("J63").Value = 999
$ws.Range("K63").Value = 1103.4
$ws.Range("L63").Value = 999
$ws.Range("M63").Value = -417.4000000000001
$ws.Range("N63").Value = -2371
$ws.Range("H66").Value = 1073.5714
$ws.Range("I66").Value = 1103.4
$ws.Range("J66").Value = 999
$ws.Range("K66").Value = 5517
$ws.Range("L66").Value = 4995
$ws.Range("M66").Value = -2085
$ws.Range("N66").Value = -11859
$ws.Range("H74").Value = 1309.5349
$ws.Range("I74").Value = 1187.1316
$ws.Range("J74").Value = 2239.8
$ws.Range("K74").Value = 1187.1316
$ws.Range("L74").Value = 2239.8
$ws.Range("M74").Value = -313.1315999999999
$ws.Range("N74").Value = -3987.8
$ws.Range("H77").Value = 1309.5349
$ws.Range("I77").Value = 1187.1316
$ws.Range("J77").Value = 2239.8
$ws.Range("K77").Value = 5935.657999999999
$ws.Range("L77").Value = 11199
$ws.Range("M77").Value = -1567.657999999999
$ws.Range("N77").Value = -19935
$ws.Range("H122").Value = 381928.5
$ws.Range("I122").Value = 3011.0938
$ws.Range("J122").Value = 2806999.8
$ws.Range("K122").Value = 9033.2814
$ws.Range("L122").Value = 8420999.399999999
$ws.Range("M122").Value = -6583.2814
$ws.Range("N122").Value = -8425899.399999999
$ws.Range("H136").Value = 11353.094
$ws.Range("I136").Value = 13332.167
$ws.Range("K136").Value = 39996.501
$ws.Range("M136").Value = -37446.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1008209.8
$ws.Range("I5").Value = 524.5
$ws.Range("J5").Value = 1680000
$ws.Range("K5").Value = 524.5
$ws.Range("L5").Value = 1680000
$ws.Range("M5").Value = -411.5
$ws.Range("N5").Value = -1680226
$ws.Range("H107").Value = 2819.6667
$ws.Range("J107").Value = 1008.6667
$ws.Range("L107").Value = 1008.6667
$ws.Range("N107").Value = -4848.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1572.8422
$ws.Range("I16").Value = 1759
$ws.Range("J16").Value = 1051.6
$ws.Range("K16").Value = 1759
$ws.Range("L16").Value = 1051.6
$ws.Range("M16").Value = -1472
$ws.Range("N16").Value = -1625.6
$ws.Range("H31").Value = 3079.6445
$ws.Range("I31").Value = 2161
$ws.Range("J31").Value = 3248.8684
$ws.Range("K31").Value = 2161
$ws.Range("L31").Value = 3248.8684
$ws.Range("M31").Value = -1866
$ws.Range("N31").Value = -3838.8684
$ws.Range("H34").Value = 3079.6445
$ws.Range("I34").Value = 2161
$ws.Range("J34").Value = 3248.8684
$ws.Range("K34").Value = 2161
$ws.Range("L34").Value = 3248.8684
$ws.Range("M34").Value = -1959
$ws.Range("N34").Value = -3652.8684
$ws.Range("H94").Value = 2209.7646
$ws.Range("I94").Value = 3581.7144
$ws.Range("J94").Value = 1249.4
$ws.Range("K94").Value = 3581.7144
$ws.Range("L94").Value = 1249.4
$ws.Range("M94").Value = -3130.7144
$ws.Range("N94").Value = -2151.4
$ws.Range("H113").Value = 1572.8422
$ws.Range("I113").Value = 1759
$ws.Range("J113").Value = 1051.6
$ws.Range("K113").Value = 1759
$ws.Range("L113").Value = 1051.6
$ws.Range("M113").Value = 411
$ws.Range("N113").Value = -5391.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 194.5
$ws.Range("I2").Value = 129.15384
$ws.Range("J2").Value = 288.8889
$ws.Range("K2").Value = 774.92304
$ws.Range("L2").Value = 1733.3334
$ws.Range("M2").Value = -661.92304
$ws.Range("N2").Value = -1959.3334
$ws.Range("H14").Value = 797.36365
$ws.Range("I14").Value = 797.36365
$ws.Range("K14").Value = 2392.09095
$ws.Range("M14").Value = -2219.09095
$ws.Range("H26").Value = 393.4
$ws.Range("I26").Value = 20
$ws.Range("J26").Value = 553.4286
$ws.Range("K26").Value = 60
$ws.Range("L26").Value = 1660.2858
$ws.Range("M26").Value = 228
$ws.Range("N26").Value = -2236.2858
$ws.Range("H38").Value = 2189.3333
$ws.Range("J38").Value = 3539.2856
$ws.Range("L38").Value = 10617.8568
$ws.Range("N38").Value = -11311.8568
$ws.Range("H92").Value = 1130.7778
$ws.Range("I92").Value = 920.25
$ws.Range("J92").Value = 1299.2
$ws.Range("K92").Value = 2760.75
$ws.Range("L92").Value = 3897.6
$ws.Range("M92").Value = -1512.75
$ws.Range("N92").Value = -6393.6
$ws.Range("H121").Value = 1053244
$ws.Range("I121").Value = 248.25
$ws.Range("K121").Value = 744.75
$ws.Range("M121").Value = 565.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6894.6772
$ws.Range("I70").Value = 6915.8887
$ws.Range("J70").Value = 6751.5
$ws.Range("K70").Value = 6915.8887
$ws.Range("L70").Value = 6751.5
$ws.Range("M70").Value = -6645.8887
$ws.Range("N70").Value = -7291.5
$ws.Range("H73").Value = 6894.6772
$ws.Range("I73").Value = 6915.8887
$ws.Range("J73").Value = 6751.5
$ws.Range("K73").Value = 6915.8887
$ws.Range("L73").Value = 6751.5
$ws.Range("M73").Value = -5979.8887
$ws.Range("N73").Value = -8623.5
$ws.Range("H107").Value = 467.33334
$ws.Range("I107").Value = 541.5
$ws.Range("J107").Value = 408
$ws.Range("K107").Value = 541.5
$ws.Range("L107").Value = 408
$ws.Range("M107").Value = 1378.5
$ws.Range("N107").Value = -4248
$ws.Range("H122").Value = 6062.0645
$ws.Range("I122").Value = 6208.069
$ws.Range("K122").Value = 18624.207
$ws.Range("M122").Value = -16174.207

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10197.611
$ws.Range("I132").Value = 15374.05
$ws.Range("J132").Value = 3727.0625
$ws.Range("K132").Value = 46122.14999999999
$ws.Range("L132").Value = 11181.1875
$ws.Range("M132").Value = -43592.14999999999
$ws.Range("N132").Value = -16241.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 33706
$ws.Range("I75").Value = 32118
$ws.Range("J75").Value = 34500
$ws.Range("K75").Value = 32118
$ws.Range("L75").Value = 34500
$ws.Range("M75").Value = -31182
$ws.Range("N75").Value = -36372
$ws.Range("H78").Value = 33706
$ws.Range("I78").Value = 32118
$ws.Range("J78").Value = 34500
$ws.Range("K78").Value = 96354
$ws.Range("L78").Value = 103500
$ws.Range("M78").Value = -91674
$ws.Range("N78").Value = -112860
$ws.Range("H132").Value = 9687.232
$ws.Range("I132").Value = 9453.312
$ws.Range("K132").Value = 28359.936
$ws.Range("M132").Value = -25829.936
